$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New column header for the Menu / non-member-supplement amounts.
$ws.Range("C1").NumberFormat = "0.00"
$ws.Range("C1").Value = "Non Member Supplement"

# Amounts for each event row (0.00 style numeric format).
$ws.Range("C2").NumberFormat = "0.00"
$ws.Range("C2").Value = 2
$ws.Range("C3").NumberFormat = "0.00"
$ws.Range("C3").Value = 1.5
$ws.Range("C4").NumberFormat = "0.00"
$ws.Range("C4").Value = 3.5

# Keep the date column using the existing built-in date format.
$ws.Range("B2:B4").NumberFormat = "d-mmm-yy"

# Widen column C to fit the new header/values.
$ws.Columns("C").ColumnWidth = 23.6

# Move the active selection to reflect the next empty data row.
$ws.Range("C5").Select() | Out-Null
